$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second data row (client 40289 / cuotas 12 / propuesta 4873666 /
# "Se han encontrado errores..." result) is removed entirely, while the
# row further down (row 7, the stray styled cell N7) keeps its original
# row number. Using Clear() on just that row's range drops the row node
# itself instead of shifting the rows below it upward.
$ws.Range("A3:V3").Clear()

# The surviving proposal's result row now reports the newly generated
# "Numero Propuesta" value and keeps it as text (it must stay a shared
# string / text cell, not be reinterpreted as a number), picking up the
# text number format/style already used elsewhere in the sheet.
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "4873680"

# Reflect the refreshed view: scrolled right so column P leads the
# visible area, with V2 (the "Resultado" cell of the remaining row) as
# the active selection.
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("V2").Select()
